$wb = $excel.ActiveWorkbook

# Add the new "Attribute" worksheet as the last sheet in the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Attribute"

# Header row
$ws.Range("A1").Value = "class"
$ws.Range("B1").Value = "name"
$ws.Range("A1:B1").Font.Bold = $true

# Data rows
$ws.Range("A2").Value = "struct"
$ws.Range("B2").Value = "Attribute"

$ws.Range("A3").Value = "field"
$ws.Range("B3").Value = "weight"

$ws.Range("B4").Value = "unit"
$ws.Range("A4").Value = "attribute"

# Match the view state used by the other sheets: freeze the header row
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A5:XFD6").Select() | Out-Null

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
